$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Extend the blank-row template (rows 18/19 currently hold data; we need
# 4 blank rows 18-21 plus a new data row at 22). Clone the already-blank,
# correctly-styled row 17 down into the two new rows 20/21 so the styles
# line up with the rest of the blank block.
$ws.Range("A17:F17").Copy($ws.Range("A20:F20"))
$ws.Range("A17:F17").Copy($ws.Range("A21:F21"))

# --- Build the new data row in row 19 (it already carries the user row's
# style s="1"/s="2"), overwriting the old YAHYA/EA4C7814/110/TRUE values with
# the new ones. For values that Excel would otherwise auto-coerce to a
# number or boolean (so they'd lose their shared-string "text" typing), set
# them as a ="..." formula and immediately flatten it to a literal value via
# copy / paste-special-values, which preserves plain text typing.
$ws.Cells.Item(19, 1).Value2 = "YAHYA"

$c = $ws.Cells.Item(19, 2)
$c.Formula = '="123"'
$c.Copy()
$c.PasteSpecial(-4163)  # xlPasteValues

$ws.Cells.Item(19, 3).Value2 = "591EF2D4"
$ws.Cells.Item(19, 4).Value2 = 0
$ws.Cells.Item(19, 5).Value2 = "103|100"

$c = $ws.Cells.Item(19, 6)
$c.Formula = '="TRUE"'
$c.Copy()
$c.PasteSpecial(-4163)  # xlPasteValues

# --- Drop the old TEST/EA4C7814/100/TRUE row.
$ws.Range("A18:F18").ClearContents()

# --- Move the newly-built row down to row 22, leaving rows 18-21 blank.
$ws.Range("A19:F19").Cut($ws.Range("A22:F22"))

$excel.CutCopyMode = 0
